$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: update title and link
$ws.Range("D12").Value = "“파이썬 라이브러리를 활용한 머신러닝(번역개정2판)”이 출간되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/02/24/%ed%8c%8c%ec%9d%b4%ec%8d%ac-%eb%9d%bc%ec%9d%b4%eb%b8%8c%eb%9f%ac%eb%a6%ac%eb%a5%bc-%ed%99%9c%ec%9a%a9%ed%95%9c-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d%eb%b2%88%ec%97%ad%ea%b0%9c%ec%a0%952%ed%8c%90/"

# Row 29: update title and link
$ws.Range("D29").Value = "[Paper Review] UNet++와 U-Net3+ 논문 리뷰"
$ws.Range("E29").Value = "https://blog.promedius.ai/unet-unet3/"
